# TCOtoliths.xlsx update:
#  1. Clear stray "Weighed?/Ran?/Note" marks (cols K, L, M) that were
#     mistakenly left on a number of existing rows.
#  2. Append five new otolith records (rows 323-327) for Macroparalepis
#     affinis caught in trawl TC6.
#  3. Move the saved selection/scroll position to the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear erroneous cells -------------------------------------------------
$cellsToClear = @(
  "K2","L2",
  "M3",
  "K4",
  "K5",
  "K6",
  "K7",
  "I133","J133","K133",
  "I136","J136","K136",
  "I138","J138","K138",
  "I140","J140","K140",
  "K147",
  "K148",
  "K153",
  "K154",
  "K155",
  "K156",
  "I158","J158","K158",
  "K161",
  "K165",
  "M166",
  "M198",
  "M199",
  "K205","L205",
  "K206",
  "M208",
  "K209",
  "K210","M210",
  "M211",
  "K214",
  "K217",
  "K218",
  "K221",
  "K288",
  "K289",
  "M292",
  "K295",
  "K297",
  "K298",
  "M299",
  "K301",
  "K303",
  "K305",
  "K306"
)

foreach ($addr in $cellsToClear) {
  $ws.Range($addr).ClearContents()
}

# --- 2. Append new records for TC6 / Macroparalepis affinis -------------------
# Column order below mirrors how the rows were actually populated in the
# source workbook (trawl/plate No. filled down first, then the sample code,
# with the one odd length value corrected last) so new shared-string
# entries land in the same order as the authored file.
$newRows = @(
  @{ Row=323; Length="76mm SL";  Mass=0.36; Well="A1"; Code="TCMA021" },
  @{ Row=324; Length="76mm SL";  Mass=0.46; Well="A2"; Code="TCMA022" },
  @{ Row=325; Length="78mm SL";  Mass=0.52; Well="A3"; Code="TCMA023" },
  @{ Row=326; Length="80mm SL";  Mass=0.46; Well="A4"; Code="TCMA024" },
  @{ Row=327; Length="122mm SL"; Mass=2.06; Well="A5"; Code="TCMA025" }
)

foreach ($r in $newRows) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value = 2
  $bCell = $ws.Cells.Item($row, 2)
  $bCell.Value = "Macroparalepis affinis"
  $bCell.Font.Italic = $true
  $ws.Cells.Item($row, 4).Value = $r.Mass
  $ws.Cells.Item($row, 5).Value = 2
  $ws.Cells.Item($row, 6).Value = $r.Well
}
foreach ($r in $newRows) {
  $ws.Cells.Item($r.Row, 7).Value = "TC6"
}
foreach ($r in $newRows) {
  $ws.Cells.Item($r.Row, 8).Value = $r.Code
}
foreach ($r in $newRows) {
  $ws.Cells.Item($r.Row, 3).Value = $r.Length
}

# --- 3. Update saved selection / scroll position -------------------------------
$null = $ws.Range("A312").Select()
try {
  $excel.ActiveWindow.ScrollRow = 312
  $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$null = $ws.Range("E328").Select()
